$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update row 8 (last data row) with refreshed metrics
$ws.Range("C8").Value = 1346
$ws.Range("D8").Value = 211
$ws.Range("E8").Value = 1135
$ws.Range("F8").Value = 8.654634946677604
$ws.Range("G8").Value = 84.32392273402675
$ws.Range("H8").Value = 15.67607726597325
